$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 3 (asset 400943320700103000 / PLC_ST_2 / 5633) is removed entirely.
# ------------------------------------------------------------------
$ws.Rows.Item(3).Delete()

# ------------------------------------------------------------------
# Row 2 data updates.
# ------------------------------------------------------------------
# B2 and C2 keep their existing "quote-prefixed text" style but the
# value becomes "000" for both (leading zeros preserved as text).
$ws.Range("B2").Value = "'000"
$ws.Range("C2").Value = "'000"

# A2 gets a brand new asset id, stored as text (quote-prefixed, large
# numeric-looking string).
$ws.Range("A2").Value = "'4009433100000000000000"

# E2 (NODE_ID) gets a new value.
$ws.Range("E2").Value = "ns=2;s=PC10.test blink"

# F2 (INDEX_LIST_VALUE) no longer holds data for this row.
$ws.Range("F2").ClearContents()

# ------------------------------------------------------------------
# Number formatting: column A (ASSET_ID) is now formatted as Text.
# Apply to A2 first, then A1, so the generated style indexes land in
# the same order as the target workbook (quote-prefixed text xf
# before the plain text xf).
# ------------------------------------------------------------------
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A1").NumberFormat = "@"

# ------------------------------------------------------------------
# Leftover column width formatting for columns E and F (kept even
# though F2 no longer has data).
# ------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 46.7
$ws.Columns.Item(6).ColumnWidth = 23.7

# ------------------------------------------------------------------
# Selection moved to F12.
# ------------------------------------------------------------------
$ws.Range("F12").Select()
